# Update "想去人数" (want-to-go count) values in the F column across sheets,
# matching the output regenerated by the gh-pages build at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 753
$ws1.Range("F4").Value = 1500
$ws1.Range("F7").Value = 147
$ws1.Range("F8").Value = 6247
$ws1.Range("F12").Value = 5212
$ws1.Range("F15").Value = 1185
$ws1.Range("F16").Value = 1185
$ws1.Range("F17").Value = 57
$ws1.Range("F18").Value = 364
$ws1.Range("F19").Value = 67
$ws1.Range("F20").Value = 11
$ws1.Range("F21").Value = 299
$ws1.Range("F22").Value = 28
$ws1.Range("F23").Value = 3716

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 83

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 83
$ws4.Range("F4").Value = 753
$ws4.Range("F5").Value = 1500
$ws4.Range("F8").Value = 147
$ws4.Range("F9").Value = 6247
$ws4.Range("F13").Value = 5212
$ws4.Range("F16").Value = 1185
$ws4.Range("F17").Value = 1185
$ws4.Range("F18").Value = 57
$ws4.Range("F19").Value = 364
$ws4.Range("F20").Value = 67
$ws4.Range("F21").Value = 11
$ws4.Range("F22").Value = 299
$ws4.Range("F23").Value = 28
$ws4.Range("F24").Value = 3716
